$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet / tab title to reflect new "through" date
$ws.Name = "Through 2022-12-21"

# Update the label in column A for the December row
$ws.Range("A13").Value = "December (through 12-21)"

# Update December row values (row 13)
$ws.Range("B13").Value = 28
$ws.Range("C13").Value = 66
$ws.Range("D13").Value = 81
$ws.Range("E13").Value = 47
$ws.Range("F13").Value = 39
$ws.Range("G13").Value = 98
$ws.Range("H13").Value = 142
$ws.Range("I13").Value = 91

# Update Total row values (row 14)
$ws.Range("B14").Value = 319
$ws.Range("C14").Value = 629
$ws.Range("D14").Value = 902
$ws.Range("E14").Value = 729
$ws.Range("F14").Value = 573
$ws.Range("G14").Value = 1362
$ws.Range("H14").Value = 1785
$ws.Range("I14").Value = 1608
